# Insert a new weekly record at row 60 ("Española" / "Segunda" from
# "Provincia de Limarí", dated 2022-07-12) and push the existing
# rows 60-77 down to 61-78, matching the new weekly report row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 60:77 down to 61:78 by inserting a new blank row at 60.
$ws.Rows("60:60").Insert()

# Populate the newly inserted row 60 with the new weekly record.
$ws.Range("A60").Value = 5
$ws.Range("B60").Value = 'Macroferia Regional de Talca'
$ws.Range("C60").Value = 'Maule'
$ws.Range("D60").Value = 44754
$ws.Range("E60").Value = 7
$ws.Range("F60").Value = 100112013
$ws.Range("G60").Value = 'Alcachofa'
$ws.Range("H60").Value = 'Española'
$ws.Range("I60").Value = 'Segunda'
$ws.Range("J60").Value = 300
$ws.Range("K60").Value = 15000
$ws.Range("L60").Value = 15000
$ws.Range("M60").Value = 15000
$ws.Range("N60").Value = '$/caja 40 unidades'
$ws.Range("O60").Value = 'Provincia de Limarí'
$ws.Range("P60").Value = 375
$ws.Range("Q60").Value = 40
$ws.Range("R60").Value = 'Hortaliza'
